$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.689860582351685
$ws.Range("B1").Value = 1.569661736488342
$ws.Range("C1").Value = 7.306740760803223
$ws.Range("D1").Value = 1.613357067108154
$ws.Range("E1").Value = 0.4718554615974426
